$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "297.19"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.92%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.35"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.35%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.010"
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.57%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07508"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.98%"
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.59%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9263"
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.91%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.401"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.17%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1220"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.63%"
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "7.17%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08840"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.27%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04112"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.49%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1054"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.13%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001282"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.48%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005777"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.98%"
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.54%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.364"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.13%"
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3317"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.75%"
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.982"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.40%"
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1419"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "5.13%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2967"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.87%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04052"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.25%"
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.20%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003879"
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.75%"
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-4.00%"
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.04%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02415"
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "4.56%"
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05213"
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.00%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.005985"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-7.13%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007806"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.68%"
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.35%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007371"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.17%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008118"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "6.44%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.2962"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-5.28%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006276"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.48%"
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.01%"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.04552"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-81.14%"
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004201"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.01%"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.01%"
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.01%"
$ws.Range("E51").Style = "Normal"
